$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 1126.324221750463
$ws.Range("E6").Value = 818.1949278344663
$ws.Range("E7").Value = 1578.747577773758
$ws.Range("E8").Value = 793.8889152819173
$ws.Range("E10").Value = 782.5306459973659
$ws.Range("E11").Value = 1585.390919286713
$ws.Range("E12").Value = 795.7781330605015
$ws.Range("E14").Value = 259.5591280471726
$ws.Range("E16").Value = 558.7753921986449
$ws.Range("E18").Value = 807.206246206097
$ws.Range("E20").Value = 784.4265302827009
$ws.Range("E22").Value = 774.1766537979651
$ws.Range("E23").Value = 1559.08339359606
$ws.Range("E24").Value = 786.4358660519543
$ws.Range("E30").Value = 830.1151962040929
$ws.Range("E32").Value = 804.3952658385263
$ws.Range("E36").Value = 806.5705040515523
$ws.Range("E38").Value = 266.5252418760376
$ws.Range("E39").Value = 584.552167130056
$ws.Range("E40").Value = 567.9523511419387
$ws.Range("E41").Value = 1126.216089402347
$ws.Range("E42").Value = 818.1793526530962
$ws.Range("E43").Value = 1578.619524922479
$ws.Range("E44").Value = 793.8847613405588
$ws.Range("E45").Value = 1580.36825256898
$ws.Range("E46").Value = 782.5229898467379
$ws.Range("E47").Value = 1585.259808063967
$ws.Range("E48").Value = 795.7744375164431
$ws.Range("E49").Value = 1588.548971149493
$ws.Range("E50").Value = 266.493025461468
$ws.Range("E52").Value = 568.0104116955799
$ws.Range("E54").Value = 813.4568609152609
$ws.Range("E55").Value = 1576.730847417499
$ws.Range("E56").Value = 793.5419290246531
$ws.Range("E57").Value = 1576.557761611099
$ws.Range("E58").Value = 782.0719234318537
$ws.Range("E59").Value = 1581.448554684218
$ws.Range("E60").Value = 795.3342777449046
$ws.Range("E62").Value = 266.3889820864327
$ws.Range("E64").Value = 566.5281526569033
$ws.Range("E66").Value = 818.0123503488376
$ws.Range("E68").Value = 797.7061678009497
$ws.Range("E70").Value = 791.5612899806395
$ws.Range("E74").Value = 261.6127993072423
$ws.Range("E76").Value = 555.7952398488034
$ws.Range("E78").Value = 798.0944966132887
$ws.Range("E80").Value = 775.541836328673
$ws.Range("E82").Value = 767.3292755626042
$ws.Range("E84").Value = 782.4442353174898
$ws.Range("E86").Value = 272.4147849066733
$ws.Range("E87").Value = 644.9778900760817
$ws.Range("E88").Value = 574.0673792294839
$ws.Range("E89").Value = 1264.585175950039
$ws.Range("E90").Value = 841.2972881187146
$ws.Range("E91").Value = 1728.55612864316
$ws.Range("E92").Value = 1034.035194554007
$ws.Range("E94").Value = 813.108517856671
$ws.Range("E95").Value = 1733.603986354046
$ws.Range("E96").Value = 828.9730477921889
$ws.Range("E97").Value = 1739.110052210355
$ws.Range("E98").Value = 266.6157337158561
$ws.Range("E100").Value = 567.102988630201
$ws.Range("E102").Value = 808.3267200468705
$ws.Range("E103").Value = 1574.902498208721
$ws.Range("E104").Value = 785.4098339593273
$ws.Range("E105").Value = 1576.594940625636
$ws.Range("E106").Value = 910.1570174712725
$ws.Range("E107").Value = 1586.070034305075
$ws.Range("E108").Value = 927.6537306076445
$ws.Range("E109").Value = 1588.836127551937
$ws.Range("E110").Value = 266.9437071860679
$ws.Range("E112").Value = 564.6747183214767
$ws.Range("E113").Value = 1247.991080168495
$ws.Range("E114").Value = 818.2605260362086
$ws.Range("E115").Value = 1566.603472635824
$ws.Range("E116").Value = 793.5246273194244
$ws.Range("E117").Value = 1568.299161055901
$ws.Range("E118").Value = 782.1971864990984
$ws.Range("E119").Value = 1573.037710344948
$ws.Range("E120").Value = 795.5405910393119
$ws.Range("E121").Value = 1576.190679413405
$ws.Range("E122").Value = 266.5276723681793
$ws.Range("E123").Value = 584.5536604936384
$ws.Range("E124").Value = 567.7037498232255
$ws.Range("E126").Value = 817.9912438009436
$ws.Range("E127").Value = 1572.195023406674
$ws.Range("E128").Value = 793.3354613206874
$ws.Range("E130").Value = 783.5242883693231
$ws.Range("E131").Value = 1674.850640197048
$ws.Range("E132").Value = 796.7313223230198
$ws.Range("E134").Value = 233.0783992990185
$ws.Range("E138").Value = 724.0398665406584
$ws.Range("E139").Value = 1434.683157889212
$ws.Range("E142").Value = 688.3316671916407
$ws.Range("E143").Value = 1440.372515638713
$ws.Range("E144").Value = 700.2573101176421
$ws.Range("E145").Value = 1443.573517806858
$ws.Range("E146").Value = 233.4162861885439
$ws.Range("E148").Value = 577.6098722051493
$ws.Range("E150").Value = 1064.661654731995
$ws.Range("E152").Value = 1061.49756278927
$ws.Range("E154").Value = 1050.469062964313
$ws.Range("E156").Value = 1051.545880453864
$ws.Range("E158").Value = 443.5581633122617
$ws.Range("E159").Value = 732.2695648967738
$ws.Range("E160").Value = 793.774860127618
$ws.Range("E161").Value = 1387.87508270595
$ws.Range("E162").Value = 852.0323127596539
$ws.Range("E163").Value = 1811.359407020379
$ws.Range("E164").Value = 830.5105738161905
$ws.Range("E165").Value = 1813.143494101459
$ws.Range("E166").Value = 818.1074879816121
$ws.Range("E167").Value = 1817.743937572714
$ws.Range("E168").Value = 827.7742181176412
$ws.Range("E169").Value = 1821.554632609486
$ws.Range("E170").Value = 255.9702156987711
$ws.Range("E171").Value = 590.3515712426746
$ws.Range("E172").Value = 549.3694482995126
$ws.Range("E174").Value = 790.946276579613
$ws.Range("E176").Value = 771.3617748906128
$ws.Range("E177").Value = 1517.870728559662
$ws.Range("E178").Value = 764.0200175228408
$ws.Range("E179").Value = 1521.748491809717
$ws.Range("E180").Value = 778.3968694238481
$ws.Range("E184").Value = 560.7005766872624
$ws.Range("E186").Value = 806.0270739693466
$ws.Range("E188").Value = 783.0044954324231
$ws.Range("E190").Value = 772.3046625696614
$ws.Range("E191").Value = 1569.275547246882
$ws.Range("E192").Value = 782.0710500932455
